$wb = $excel.ActiveWorkbook

# Rename the commodity label "Copper ores and concentrates" -> "Copper"
# (disaggregation of the Copper commodity) everywhere it appears, i.e. in
# cell C4 of every year sheet.
foreach ($ws in $wb.Worksheets) {
    $c4 = $ws.Range("C4")
    if ($c4.Value2 -eq "Copper ores and concentrates") {
        $c4.Value = "Copper"
    }
}

# A handful of year-sheets also carry a re-computed Copper figure (D4) that
# shifted by a tiny (last-significant-digit) amount as part of the same
# disaggregation re-export.
$sheet2022 = $wb.Worksheets.Item("2022")
$sheet2022.Range("D4").Value = 13486.09352961935

$sheet2039 = $wb.Worksheets.Item("2039")
$sheet2039.Range("D4").Value = 240901.536595003

$sheet2067 = $wb.Worksheets.Item("2067")
$sheet2067.Range("D4").Value = 831477.5295182781

$sheet2075 = $wb.Worksheets.Item("2075")
$sheet2075.Range("D4").Value = 1847958.449639017

$sheet2082 = $wb.Worksheets.Item("2082")
$sheet2082.Range("D4").Value = 1706494.216939882
